# Fruta / hortaliza, semanal
# Inserts two new daily price rows for "Feria Lagunitas de Puerto Montt - Plátano"
# right after the existing row 744, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 745:746 - everything currently at row 745 downward
# shifts down by two rows (old 745 -> new 747, ..., old 853 -> new 855).
$ws.Rows("745:746").Insert()

# Populate the first newly inserted row (745)
$ws.Cells.Item(745, 1).Value = 4
$ws.Cells.Item(745, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(745, 3).Value = "Los Lagos"
$ws.Cells.Item(745, 4).Value = 45077
$ws.Cells.Item(745, 5).Value = 10
$ws.Cells.Item(745, 6).Value = "Fruta"
$ws.Cells.Item(745, 7).Value = 100108
$ws.Cells.Item(745, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(745, 9).Value = 100108006
$ws.Cells.Item(745, 10).Value = "Plátano"
$ws.Cells.Item(745, 11).Value = "Sin especificar"
$ws.Cells.Item(745, 12).Value = "Pintón"
$ws.Cells.Item(745, 13).Value = 100
$ws.Cells.Item(745, 14).Value = 19000
$ws.Cells.Item(745, 15).Value = 19000
$ws.Cells.Item(745, 16).Value = 19000
$ws.Cells.Item(745, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(745, 18).Value = "Ecuador"
$ws.Cells.Item(745, 19).Value = 950
$ws.Cells.Item(745, 20).Value = 20

# Populate the second newly inserted row (746)
$ws.Cells.Item(746, 1).Value = 4
$ws.Cells.Item(746, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(746, 3).Value = "Los Lagos"
$ws.Cells.Item(746, 4).Value = 45077
$ws.Cells.Item(746, 5).Value = 10
$ws.Cells.Item(746, 6).Value = "Fruta"
$ws.Cells.Item(746, 7).Value = 100108
$ws.Cells.Item(746, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(746, 9).Value = 100108006
$ws.Cells.Item(746, 10).Value = "Plátano"
$ws.Cells.Item(746, 11).Value = "Sin especificar"
$ws.Cells.Item(746, 12).Value = "Primera Pintón"
$ws.Cells.Item(746, 13).Value = 200
$ws.Cells.Item(746, 14).Value = 20000
$ws.Cells.Item(746, 15).Value = 21000
$ws.Cells.Item(746, 16).Value = 20500
$ws.Cells.Item(746, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(746, 18).Value = "Ecuador"
$ws.Cells.Item(746, 19).Value = 1025
$ws.Cells.Item(746, 20).Value = 20

Write-Output "rows inserted and populated"
